# Update the "error margin" tables to use a 99% confidence interval (z=2.58)
# computed over 100k runs instead of a 95% interval (z=1.96) over 1m runs.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Required Tasks")
$ws2 = $wb.Worksheets.Item("Desirable Tasks")

# ---------------------------------------------------------------------------
# Sheet "Required Tasks": BasicWeightTable
# ---------------------------------------------------------------------------

$ws1.Range("D2").Value = "99% confidence level"
$ws1.Range("D3").Value = "100k runs error margin"

$ws1.Range("F4").Formula = "=SQRT(100000 * C4 * (1-C4)) * 2.58"
$ws1.Range("F5:F7").Formula = "=SQRT(100000 * C5 * (1-C5)) * 2.58"
$ws1.Range("G4").Formula = "=F4/100000"
$ws1.Range("G5:G7").Formula = "=F5/100000"

$ws1.Range("D4").Value = 0.003996918813286054
$ws1.Range("D5").Value = 0.0037387752005168754
$ws1.Range("D6").Value = 0.0032634705452937677
$ws1.Range("D7").Value = 0.0024476029089703253
$ws1.Range("D4:D7").NumberFormat = "0.0000%"

# ---------------------------------------------------------------------------
# Sheet "Desirable Tasks": WinRules
# ---------------------------------------------------------------------------

$ws2.Range("L2").Value = "99% confidence level"
$ws2.Range("L3").Value = "100k runs error margin"

$ws2.Range("K11").Formula = "=SQRT(100000*K4*(1-K4))*2.58"
$ws2.Range("K12:K15").Formula = "=SQRT(100000*K5*(1-K5))*2.58"
$ws2.Range("L11").Formula = "=K11/100000"
$ws2.Range("L12:L15").Formula = "=K12/100000"

$ws2.Range("L4").Value = 0.0005538411045921045
$ws2.Range("L5").Value = 0.0008499723604911371
$ws2.Range("L6").Value = 0.0011815745644745545
$ws2.Range("L7").Value = 0.0015407870363415006
$ws2.Range("L8").Value = 0.0020740041117499153

$wb.Save()
